$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 89, shifting rows 89:173 down to 90:174
$ws.Rows.Item(89).Insert()

# The new row 89 receives what used to be row 88's data (the row above it)
$ws.Range("A88:R88").Copy()
$ws.Range("A89").PasteSpecial()

# Update row 88 (which keeps its row number) with the new date and volume
$ws.Cells.Item(88, 4).Value = 44586
$ws.Cells.Item(88, 10).Value = 55
